$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." built from the runs:
#   "Versi" | "on" | " 2" | "." (with bookmark _GoBack between " 2" and ".")
# Target: "Version 1." built from the runs:
#   "Version" | " 1." (bookmark _GoBack kept after " 1.", trailing "." run removed)

# Step 1: merge "Versi" + "on" into a single "Version" run.
$r1 = $d.Range(5, 7)
$r1.Delete()
$r2 = $d.Range(0, 5)
$r2.InsertAfter("on")

# Step 2: change the "2" digit run into "1." (keeps leading space).
$r3 = $d.Range(7, 9)
$r3.Text = " 1."

# Step 3: remove the now-redundant trailing "." run.
$r4 = $d.Range(10, 11)
$r4.Delete()
